# Update brake actuator details (row 8) on the "power" sheet, and add a
# new "BOTS" switch entry, per commit "update brake actuator details".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: brake motor details -------------------------------------------------
# Device/part name (column B) and note (column C) for the brake motor.
$ws.Range("B8").Value = "Faulhaber DC-Micromotors`nSeries 3890H024CR"
$ws.Range("C8").Value = " Powered from steering contoller. Planetary Gearheads`nSeries 38/2 S"

# Apply the highlighted (yellow), centered, wrapped style to B8:C8.
$fmtRange = $ws.Range("B8:C8")
$fmtRange.Interior.Color = 65535
$fmtRange.HorizontalAlignment = -4108
$fmtRange.VerticalAlignment = -4108
$fmtRange.WrapText = $true

# Updated max wattage for the brake motor.
$ws.Range("F8").Value = 160

# Row height grows to fit the new two-line text.
$ws.Rows.Item(8).RowHeight = 46.5

# --- New switch entry: BOTS ------------------------------------------------------
$ws.Range("A34").Value = "BOTS"

# --- Sheet view: zoom + selection -------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("C3").Select() | Out-Null
$excel.ActiveWindow.Zoom = 90
